# Remove the containing path prefix "sequence/run_0731_samples/" from the
# fastq filenames stored in column F (rows 2-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "sequence/run_0731_samples/"

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 6)  # Column F
    $value = $cell.Value2
    if ($value -ne $null -and $value.ToString().StartsWith($prefix)) {
        $cell.Value2 = $value.ToString().Substring($prefix.Length)
    }
}

# Update the active selection to match the recorded state in the file.
$ws.Range("F20").Select()
